$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G11:G18 values to 51
$ws.Range("G11").Value = 51
$ws.Range("G12").Value = 51
$ws.Range("G13").Value = 51
$ws.Range("G14").Value = 51
$ws.Range("G15").Value = 51
$ws.Range("G16").Value = 51
$ws.Range("G17").Value = 51
$ws.Range("G18").Value = 51

# Update selection / active cell to G19
$ws.Range("G19").Select()
